$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Age based discount row: "20,60" -> "20,50"
$ws.Range("B11").Value = "20,50"

# Amount based discount row: 700 -> 500
$ws.Range("E14").Value = 500

# Update the active selection to B11 (matches the author's final click)
$null = $ws.Range("B11").Select()
